# Update crypto volume(1h) percentages in column E per the latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  +6.33%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  +4.85%  "
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("E37").Value = "  +3.84%  "
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +4.28%  "
